# Nayem_meal.xlsx - "meal finished at 23rd March"
# Record meals (2.5) for day 22 (column W) and day 23 (column X)
# for every person listed in rows 3-9, then leave the selection on X3:X9
# (matching where the author was last working in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("W3:X9").Value = 2.5

$ws.Range("X3:X9").Select()
